$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Part 1: "This sample is compatible with the Windows 10 Creators
# Update SDK (15063)" paragraph -> split into several italic runs with
# a _GoBack bookmark wrapping "Windows 10 Fall Creators Update SDK
# (16299)", and the number/word updated to reflect the Fall Creators
# Update (16299).
# ------------------------------------------------------------------

$oldIntro = "This sample is compatible with the Windows 10 Creators Update SDK (15063)"
$newIntro = "This sample is compatible with the Windows 10 Fall Creators Update SDK (16299)"
$found = $d.Content.Find.Execute($oldIntro, $true, $false, $false, $false, $false, $true, 1, $false, $newIntro, 2)

# Locate the paragraph again and compute absolute offsets for each run
# boundary within the replaced sentence.
$introPara = $d.Paragraphs.Item(2)
$introStart = $introPara.Range.Start

$segments = @(
    "This sample is compatible with the ",
    "Windows 10",
    " Fall",
    " Creators Update SDK (",
    "16299",
    ")"
)

$offsets = @(0)
$running = 0
foreach ($seg in $segments) {
    $running += $seg.Length
    $offsets += $running
}

# Insert zero-length temporary bookmarks at each interior boundary so
# that the run gets split into one run per segment.
for ($i = 1; $i -lt ($offsets.Length - 1); $i++) {
    $pos = $introStart + $offsets[$i]
    $marker = $d.Range($pos, $pos)
    $d.Bookmarks.Add("TempIntroSplit$i", $marker)
}

# Remove the temporary split bookmarks, leaving the runs separated.
for ($i = 1; $i -lt ($offsets.Length - 1); $i++) {
    $d.Bookmarks.Item("TempIntroSplit$i").Delete()
}

# Wrap the "Windows 10 Fall Creators Update SDK (16299)" portion (all
# segments after the first) with the _GoBack bookmark.
$bmStart = $introStart + $offsets[1]
$bmEnd = $introStart + $offsets[$offsets.Length - 1]
$bmRange = $d.Range($bmStart, $bmEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ------------------------------------------------------------------
# Part 2: Merge the five runs making up "Running this sample requires
# the Windows 10 Anniversary Update (14393) or later." into a single
# run (keeping the preceding " " run separate), and remove the old
# _GoBack bookmark paragraph, turning it into a plain empty paragraph.
# ------------------------------------------------------------------

$reqParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.IndexOf("Running this sample requires") -ge 0) {
        $reqParaIndex = $i
        break
    }
}

$reqPara = $d.Paragraphs.Item($reqParaIndex)
$reqText = $reqPara.Range.Text
$reqIdx = $reqText.IndexOf("Running this sample requires")
$reqParaStart = $reqPara.Range.Start
$reqRunStart = $reqParaStart + $reqIdx
$reqRunEnd = $reqPara.Range.End - 1

# Guard the left boundary with a temporary bookmark so the replacement
# does not get merged back into the preceding " " run.
$guardMarker = $d.Range($reqRunStart, $reqRunStart)
$d.Bookmarks.Add("TempReqGuard", $guardMarker)

# First set to a placeholder distinct from the original text so the
# engine actually performs (and records) the run merge...
$reqRange = $d.Range($reqRunStart, $reqRunEnd)
$reqRange.Text = "TempPlaceholderRunMergeText"

# ...then relocate the placeholder and set the final desired text.
$reqPara2 = $d.Paragraphs.Item($reqParaIndex)
$reqText2 = $reqPara2.Range.Text
$reqIdx2 = $reqText2.IndexOf("TempPlaceholderRunMergeText")
$reqParaStart2 = $reqPara2.Range.Start
$reqRunStart2 = $reqParaStart2 + $reqIdx2
$reqRunEnd2 = $reqPara2.Range.End - 1
$reqRange2 = $d.Range($reqRunStart2, $reqRunEnd2)
$reqRange2.Text = "Running this sample requires the Windows 10 Anniversary Update (14393) or later."

$d.Bookmarks.Item("TempReqGuard").Delete()

# Note: the old _GoBack bookmark (which lived alone in its own
# paragraph further down the document) was already implicitly removed
# above when we (re)added a bookmark named "_GoBack" around the
# "Windows 10 Fall Creators Update SDK (16299)" text - bookmark names
# are unique per document, so adding it again relocates it and leaves
# a plain empty paragraph behind at its old location.
